$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing column E header ("QADF significance (10%)" -> "Mean Reversion (Entire period)")
# and add the new column F header ("Mean Reversion (EAEU only)")
$ws.Range("E1").Value = "Mean Reversion (Entire period)"
$ws.Range("F1").Value = "Mean Reversion (EAEU only)"

# Fill in the new column F values for each country row
$ws.Range("F2").Value = "N/A"
$ws.Range("F3").Value = "N/A"
$ws.Range("F4").Value = "0.1-0.4"
$ws.Range("F5").Value = "0.7,0.8"
# "0.9" looks numeric, force it to stay text like the rest of the column
$ws.Range("F6").Value = "'0.9"
$ws.Range("F6").Style = "Normal"
$ws.Range("F7").Value = "0.5-0.9"
$ws.Range("F8").Value = "0.2-0.7"
$ws.Range("F9").Value = "N/A"

# Extend the conditional formatting ranges that previously stopped at column E
# so that they now also cover the new column F
$fcs = $ws.Range("A1:F9").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq '$A$1:$E$1') {
        $fc.ModifyAppliesToRange($ws.Range("A1:F1"))
    }
    elseif ($addr -eq '$A$1:$E$9') {
        $fc.ModifyAppliesToRange($ws.Range("A1:F9"))
    }
    elseif ($addr -eq '$A$9:$E$9') {
        $fc.ModifyAppliesToRange($ws.Range("A9:F9"))
    }
}
